# Implementation of suggestions of Biostat reviewer
# - Widen the third table column (to fit the longer, more precise values)
# - Replace the IRR (95% CI), p-value figures with more precise ones

$d = $word.ActiveDocument

# Widen the 3rd grid column from 3254 twips (162.7pt) to 3622 twips (181.1pt)
$table = $d.Tables.Item(1)
$table.Columns.Item(3).Width = 181.1

$wdReplaceOne = 1
$wdFindContinue = 1

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, $wdFindContinue, $false, $new, $wdReplaceOne) | Out-Null
}

Replace-Text "0.54 (0.46 to 0.65), p < 0.001" "0.543 (0.457 to 0.646), p < 0.001"
Replace-Text "1.04 (1.03 to 1.06), p < 0.001" "1.044 (1.03 to 1.057), p < 0.001"
Replace-Text "1.18 (1.09 to 1.28), p < 0.001" "1.182 (1.094 to 1.277), p < 0.001"
Replace-Text "1.03 (1.02 to 1.03), p < 0.001" "1.026 (1.02 to 1.031), p < 0.001"
Replace-Text "1.22 (1.18 to 1.26), p < 0.001" "1.218 (1.18 to 1.257), p < 0.001"
Replace-Text "1.01 (1.01 to 1.01), p < 0.001" "1.008 (1.006 to 1.01), p < 0.001"
Replace-Text "3.3 (0.48 to 22.65), p = 0.282" "1.206 (0.204 to 7.108), p = 0.85"
Replace-Text "1 (0.92 to 1.08), p = 0.934" "0.997 (0.891 to 1.116), p = 0.967"
Replace-Text "1.02 (0.96 to 1.09), p = 0.536" "1.117 (1.039 to 1.201), p = 0.008"
Replace-Text "1.02 (1.01 to 1.03), p = 0.003" "1.003 (0.998 to 1.008), p = 0.283"
Replace-Text "1.55 (1.37 to 1.75), p < 0.001" "1.552 (1.373 to 1.753), p < 0.001"
Replace-Text "1.01 (1 to 1.03), p = 0.109" "1.013 (0.999 to 1.026), p = 0.109"
Replace-Text "1.17 (1.02 to 1.34), p = 0.045" "1.168 (1.017 to 1.341), p = 0.045"
Replace-Text "0.98 (0.97 to 0.99), p = 0.008" "0.984 (0.974 to 0.995), p = 0.008"
Replace-Text "0.89 (0.75 to 1.06), p = 0.217" "0.891 (0.753 to 1.055), p = 0.217"
Replace-Text "1.01 (1 to 1.03), p = 0.056" "1.013 (1.001 to 1.026), p = 0.056"

Write-Host "Done applying Biostat reviewer edits"
